# Updated cryptos list (Price and Volume(1h) columns) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.887.94"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.634.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.95"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.95"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.09%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0915"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.869.76"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.636.08"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.32"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +19.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.88"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.918.19"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.31"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.70"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0702"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.93"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.44%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.80%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.10"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.60"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.110"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.62"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.45%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.04%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.19"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.430.19"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.12%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.80"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.17%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "77.10"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +16.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0171"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.555"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.00"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0493"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "53.85"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.89%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.776.77"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.35"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "89.59"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0111"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.76%  "
